# Update the group member list in both the document body and the header.
#
# Body (first paragraph): two runs
#   "Group- Evans, Oyo, Alex, Caro and " + "Rhona"
#   -> "Group- " + "Evans, Alex, Caroline"
#
# Header (red bold banner): single run
#   "Group- Evans, Oyo, Alex, Caro and Rhona"
#   -> "Group- Evans, Alex, Caro" + "line" + " "   (split to mirror the diff)

$d = $word.ActiveDocument

# --- Body text -------------------------------------------------------
$d.Content.Find.Execute(
    "Group- Evans, Oyo, Alex, Caro and ", $true, $false, $false, $false,
    $false, $true, 1, $false, "Group- ", 2)

$d.Content.Find.Execute(
    "Rhona", $true, $false, $false, $false,
    $false, $true, 1, $false, "Evans, Alex, Caroline", 2)

# --- Header text -------------------------------------------------------
$section = $d.Sections.First
$header = $section.Headers.Item(1)

$header.Range.Find.Execute(
    "Group- Evans, Oyo, Alex, Caro and Rhona", $true, $false, $false, $false,
    $false, $true, 1, $false, "Group- Evans, Alex, Caroline ", 2)
